# Refresh crypto price ("D") and 1h volume-change ("E") figures on Sheet1
# to match the latest GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.878.16"
$ws.Range("E2").Value = "  +4.31%  "
$ws.Range("D3").Value = "2.673.47"
$ws.Range("E3").Value = "  +7.77%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  +9.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "326.12"
$ws.Range("E6").Value = "  +2.97%  "
$ws.Range("E7").Value = "  +2.03%  "
$ws.Range("E9").Value = "  +3.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.84"
$ws.Range("E10").Value = "  +5.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.11"
$ws.Range("E11").Value = "  -1.27%  "
$ws.Range("E12").Value = "  +3.34%  "
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.37"
$ws.Range("E14").Value = "  +5.20%  "
$ws.Range("D15").Value = "3.093.79"
$ws.Range("E15").Value = "  +7.67%  "
$ws.Range("D16").Value = "2.692.75"
$ws.Range("E16").Value = "  +13.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.875"
$ws.Range("E17").Value = "  +6.52%  "
$ws.Range("D18").Value = "49.857.32"
$ws.Range("E18").Value = "  +4.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.16"
$ws.Range("E19").Value = "  +4.18%  "
$ws.Range("E20").Value = "  +4.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.91"
$ws.Range("E21").Value = "  +0.97%  "
$ws.Range("E22").Value = "  +3.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "277.08"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.85"
$ws.Range("E24").Value = "  +1.94%  "
$ws.Range("E25").Value = "  +3.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.83"
$ws.Range("E26").Value = "  +5.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.20"
$ws.Range("E28").Value = "  +6.83%  "
$ws.Range("E29").Value = "  +1.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.21"
$ws.Range("E30").Value = "  +5.71%  "
$ws.Range("E31").Value = "  +3.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.28"
$ws.Range("E32").Value = "  +2.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.48"
$ws.Range("E33").Value = "  +4.88%  "
$ws.Range("E34").Value = "  +3.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0808"
$ws.Range("E35").Value = "  +5.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.07"
$ws.Range("E36").Value = "  +13.30%  "
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("E38").Value = "  +7.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.13"
$ws.Range("E39").Value = "  +9.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "124.88"
$ws.Range("E40").Value = "  +3.91%  "
$ws.Range("E41").Value = "  +2.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.53"
$ws.Range("E42").Value = "  +5.74%  "
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0317"
$ws.Range("E44").Value = "  +6.93%  "
$ws.Range("D45").Value = "2.120.79"
$ws.Range("E45").Value = "  +7.11%  "
$ws.Range("E46").Value = "  +7.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.23"
$ws.Range("E47").Value = "  +8.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.05"
$ws.Range("E48").Value = "  +7.68%  "
$ws.Range("E49").Value = "  +1.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.33"
$ws.Range("E50").Value = "  +5.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "59.30"
$ws.Range("E51").Value = "  +7.18%  "
